$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 80 (ALC)
$ws.Range("H80").Value = 703.3077
$ws.Range("I80").Value = 373.9
$ws.Range("J80").Value = 909.1875
$ws.Range("K80").Value = 1121.7
$ws.Range("L80").Value = 2727.5625
$ws.Range("M80").Value = -123.6999999999998
$ws.Range("N80").Value = -4723.5625

# Row 83 (ALC)
$ws.Range("H83").Value = 703.3077
$ws.Range("I83").Value = 373.9
$ws.Range("J83").Value = 909.1875
$ws.Range("K83").Value = 3365.1
$ws.Range("L83").Value = 8182.6875
$ws.Range("M83").Value = 1626.9
$ws.Range("N83").Value = -18166.6875

# Row 86 (ALC)
$ws.Range("H86").Value = 2533
$ws.Range("I86").Value = 4700
$ws.Range("J86").Value = 1449.5
$ws.Range("K86").Value = 4700
$ws.Range("L86").Value = 1449.5
$ws.Range("M86").Value = -3577
$ws.Range("N86").Value = -3695.5

# Row 88 (ALC)
$ws.Range("H88").Value = 1744.1666
$ws.Range("I88").Value = 3654
$ws.Range("J88").Value = 1362.2
$ws.Range("K88").Value = 3654
$ws.Range("L88").Value = 1362.2
$ws.Range("M88").Value = -3248
$ws.Range("N88").Value = -2174.2

# Row 89 (ALC)
$ws.Range("H89").Value = 2533
$ws.Range("I89").Value = 4700
$ws.Range("J89").Value = 1449.5
$ws.Range("K89").Value = 23500
$ws.Range("L89").Value = 7247.5
$ws.Range("M89").Value = -17884
$ws.Range("N89").Value = -18479.5

# Row 91 (ALC)
$ws.Range("H91").Value = 1744.1666
$ws.Range("I91").Value = 3654
$ws.Range("J91").Value = 1362.2
$ws.Range("K91").Value = 3654
$ws.Range("L91").Value = 1362.2
$ws.Range("M91").Value = -2250
$ws.Range("N91").Value = -4170.2

# Row 113 (ALC)
$ws.Range("H113").Value = 1350.75
$ws.Range("I113").Value = 1350.75
$ws.Range("K113").Value = 1350.75
$ws.Range("M113").Value = 1903.25

# Row 115 (ALC)
$ws.Range("H115").Value = 634.3333
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").Value = ""

# Row 135 (ALC)
$ws.Range("H135").Value = 1920.6666
$ws.Range("I135").Value = 1988.5
$ws.Range("J135").Value = 1866.4
$ws.Range("K135").Value = 17896.5
$ws.Range("L135").Value = 16797.6
$ws.Range("M135").Value = -15361.5
$ws.Range("N135").Value = -21867.6

$ws = $wb.Worksheets.Item("ARM")
# Row 74 (ARM)
$ws.Range("H74").Value = 9285.429
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = ""

# Row 77 (ARM)
$ws.Range("H77").Value = 9285.429
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = ""

# Row 96 (ARM)
$ws.Range("H96").Value = 6695231.5
$ws.Range("J96").Value = 6695231.5
$ws.Range("L96").Value = 6695231.5
$ws.Range("N96").Value = -6700723.5

# Row 122 (ARM)
$ws.Range("H122").Value = 1611.6
$ws.Range("I122").Value = 1611.6
$ws.Range("K122").Value = 4834.799999999999
$ws.Range("M122").Value = -2384.799999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (BSM)
$ws.Range("H86").Value = 6988.3
$ws.Range("J86").Value = 9998.833000000001
$ws.Range("L86").Value = 9998.833000000001
$ws.Range("N86").Value = -12244.833

# Row 89 (BSM)
$ws.Range("H89").Value = 6988.3
$ws.Range("J89").Value = 9998.833000000001
$ws.Range("L89").Value = 49994.165
$ws.Range("N89").Value = -61226.165

# Row 102 (BSM)
$ws.Range("H102").Value = 6602.8335
$ws.Range("I102").Value = 6602.8335
$ws.Range("K102").Value = 6602.8335
$ws.Range("M102").Value = -3357.8335

# Row 107 (BSM)
$ws.Range("H107").Value = 5411.55
$ws.Range("I107").Value = 1703.6666
$ws.Range("K107").Value = 1703.6666
$ws.Range("M107").Value = 216.3334

$ws = $wb.Worksheets.Item("CRP")
# Row 107 (CRP)
$ws.Range("H107").Value = 434.6154
$ws.Range("I107").Value = 583.3333
$ws.Range("J107").Value = 307.14285
$ws.Range("K107").Value = 583.3333
$ws.Range("L107").Value = 307.14285
$ws.Range("M107").Value = 1336.6667
$ws.Range("N107").Value = -4147.14285

$ws = $wb.Worksheets.Item("CUL")
# Row 23 (CUL)
$ws.Range("H23").Value = 654.375
$ws.Range("I23").Value = 100
$ws.Range("J23").Value = 839.1667
$ws.Range("K23").Value = 300
$ws.Range("L23").Value = 2517.5001
$ws.Range("M23").Value = -65
$ws.Range("N23").Value = -2987.5001

# Row 87 (CUL)
$ws.Range("H87").Value = 249.5
$ws.Range("I87").Value = 249.5
$ws.Range("K87").Value = 748.5
$ws.Range("M87").Value = 499.5

# Row 90 (CUL)
$ws.Range("H90").Value = 249.5
$ws.Range("I90").Value = 249.5
$ws.Range("K90").Value = 2245.5
$ws.Range("M90").Value = 3994.5

# Row 128 (CUL)
$ws.Range("H128").Value = 424999.5
$ws.Range("I128").Value = 424999.5
$ws.Range("K128").Value = 1274998.5
$ws.Range("M128").Value = -1270018.5

$ws = $wb.Worksheets.Item("GSM")
# Row 36 (GSM)
$ws.Range("H36").Value = 10000
$ws.Range("J36").Value = 10000
$ws.Range("L36").Value = 10000
$ws.Range("N36").Value = -10970

# Row 70 (GSM)
$ws.Range("H70").Value = 5523
$ws.Range("I70").Value = 4619.8887
$ws.Range("J70").Value = 8232.333000000001
$ws.Range("K70").Value = 4619.8887
$ws.Range("L70").Value = 8232.333000000001
$ws.Range("M70").Value = -4349.8887
$ws.Range("N70").Value = -8772.333000000001

# Row 73 (GSM)
$ws.Range("H73").Value = 5523
$ws.Range("I73").Value = 4619.8887
$ws.Range("J73").Value = 8232.333000000001
$ws.Range("K73").Value = 4619.8887
$ws.Range("L73").Value = 8232.333000000001
$ws.Range("M73").Value = -3683.8887
$ws.Range("N73").Value = -10104.333

# Row 80 (GSM)
$ws.Range("H80").Value = 2212.5
$ws.Range("I80").Value = 1075
$ws.Range("J80").Value = 3350
$ws.Range("K80").Value = 1075
$ws.Range("L80").Value = 3350
$ws.Range("M80").Value = -77
$ws.Range("N80").Value = -5346

# Row 83 (GSM)
$ws.Range("H83").Value = 2212.5
$ws.Range("I83").Value = 1075
$ws.Range("J83").Value = 3350
$ws.Range("K83").Value = 5375
$ws.Range("L83").Value = 16750
$ws.Range("M83").Value = -383
$ws.Range("N83").Value = -26734

# Row 126 (GSM)
$ws.Range("H126").Value = 3687.5715
$ws.Range("I126").Value = 3385.5
$ws.Range("J126").Value = 5500
$ws.Range("K126").Value = 10156.5
$ws.Range("L126").Value = 16500
$ws.Range("M126").Value = -7686.5
$ws.Range("N126").Value = -21440

# Row 132 (GSM)
$ws.Range("H132").Value = 5011.9546
$ws.Range("I132").Value = 5290.364
$ws.Range("K132").Value = 15871.092
$ws.Range("M132").Value = -13341.092

$ws = $wb.Worksheets.Item("LTW")
# Row 47 (LTW)
$ws.Range("H47").Value = 40500
$ws.Range("J47").Value = 51000
$ws.Range("L47").Value = 51000
$ws.Range("N47").Value = -51980

# Row 52 (LTW)
$ws.Range("H52").Value = 40500
$ws.Range("J52").Value = 51000
$ws.Range("L52").Value = 51000
$ws.Range("N52").Value = -51466

# Row 68 (LTW)
$ws.Range("H68").Value = 9812.5
$ws.Range("I68").Value = 9500
$ws.Range("J68").Value = 9916.666999999999
$ws.Range("K68").Value = 9500
$ws.Range("L68").Value = 9916.666999999999
$ws.Range("M68").Value = -8751
$ws.Range("N68").Value = -11414.667

# Row 71 (LTW)
$ws.Range("H71").Value = 9812.5
$ws.Range("I71").Value = 9500
$ws.Range("J71").Value = 9916.666999999999
$ws.Range("K71").Value = 47500
$ws.Range("L71").Value = 49583.335
$ws.Range("M71").Value = -43756
$ws.Range("N71").Value = -57071.335

$ws = $wb.Worksheets.Item("WVR")
# Row 62 (WVR)
$ws.Range("H62").Value = 11600
$ws.Range("I62").Value = 10000
$ws.Range("K62").Value = 10000
$ws.Range("M62").Value = -9376

# Row 65 (WVR)
$ws.Range("H65").Value = 11600
$ws.Range("I65").Value = 10000
$ws.Range("K65").Value = 50000
$ws.Range("M65").Value = -46880

# Row 126 (WVR)
$ws.Range("H126").Value = 4241.8
$ws.Range("I126").Value = 1171.5
$ws.Range("K126").Value = 3514.5
$ws.Range("M126").Value = -1044.5
